$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.356.22"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.359.30"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'519.57"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "'135.66"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("E10").Value = "  +4.73%  "
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "'24.38"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "2.779.35"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "57.345.12"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "2.360.76"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "'328.93"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").Value = "'4.24"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").Value = "'6.74"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'61.32"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "'8.84"
$ws.Range("E24").Value = "  +13.27%  "
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("D26").Value = "'0.996"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "'1.36"
$ws.Range("E27").Value = "  +11.41%  "
$ws.Range("D28").Value = "0.0₃0743"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "'167.47"
$ws.Range("E29").Value = "  -2.76%  "
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "'6.29"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +2.17%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  -3.50%  "
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  +5.63%  "
$ws.Range("E39").Value = "  +3.02%  "
$ws.Range("D40").Value = "'149.78"
$ws.Range("E40").Value = "  +6.94%  "
$ws.Range("D41").Value = "'0.383"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").Value = "'5.33"
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("D44").Value = "'284.56"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").Value = "'0.0940"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").Value = "'18.28"
$ws.Range("E48").Value = "  +4.83%  "
$ws.Range("D49").Value = "'0.0219"
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'17.67"
$ws.Range("E50").Value = "  +3.45%  "
$ws.Range("B51").Value = "Polygon"
$ws.Range("C51").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D51").Value = "'0.362"
$ws.Range("E51").Value = "  -5.85%  "
